$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B27: convert from text "3" to numeric 3 (matches rest of column)
$ws.Range("B27").Value = 3

# Add new row 28
$ws.Range("A28").Value = "Sunsi Wu"

# B28 stays a text "4" (mirrors the pre-existing text-typed anomaly),
# so force text formatting before assigning, then restore default style.
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "4"
$ws.Range("B28").Style = "Normal"

$ws.Range("C28").Value = "appealing;important"
$ws.Range("D28").Value = "DIS"
$ws.Range("E28").Value = "MET"
$ws.Range("F28").Value = "f913699b-da49-47c6-8043-88c593733ae2"
$ws.Range("G28").Value = "BJyy3a0Ez_annotated.xlsx"
$ws.Range("H28").Value = "The idea of model-parallelism (as opposed to data parallelism) is appealing and an important open problem."
